$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = -0.35038893185310371
$ws.Range("B1").Value = 0.3493729848110263
$ws.Range("A2").Value = -0.22451933514038558
$ws.Range("B2").Value = 0.22187204317006781
$ws.Range("A3").Value = -0.11892048941767897
$ws.Range("B3").Value = 0.1181532464354369
$ws.Range("A4").Value = -0.18015138600264535
$ws.Range("B4").Value = 0.17894941172827217
$ws.Range("A5").Value = -0.17294941188288071
$ws.Range("B5").Value = 0.17051819492469633
$ws.Range("A6").Value = -0.069531305379319175
$ws.Range("B6").Value = 0.069463608093902174
$ws.Range("A7").Value = -0.04946360828808416
$ws.Range("B7").Value = 0.049338281854980792
$ws.Range("A8").Value = -0.068797857938643681
$ws.Range("B8").Value = 0.068395210591020827
$ws.Range("A9").Value = -0.062395210756493569
$ws.Range("B9").Value = 0.062052297970669912
$ws.Range("A10").Value = -0.05605229813897239
$ws.Range("B10").Value = 0.056000944537458963
$ws.Range("A11").Value = -0.051500944702599583
$ws.Range("B11").Value = 0.051417040309488016
$ws.Range("A12").Value = -0.045417040478946902
$ws.Range("B12").Value = 0.045160253778523618
$ws.Range("A13").Value = -0.039160253950949908
$ws.Range("B13").Value = 0.039090170664909607
$ws.Range("A14").Value = -0.027090170852485329
$ws.Range("B14").Value = 0.027056060816157412
$ws.Range("A15").Value = -0.021056060990221503
$ws.Range("B15").Value = 0.021029187641252278
$ws.Range("A16").Value = -0.015029187815982503
$ws.Range("B16").Value = 0.015004908327540978
$ws.Range("A17").Value = -0.0090049085031402853
$ws.Range("B17").Value = 0.0089999998170604556
$ws.Range("A18").Value = -0.082358698790077511
$ws.Range("B18").Value = 0.082252119222147257
$ws.Range("A19").Value = -0.073252119375805336
$ws.Range("B19").Value = 0.072439365057083371
$ws.Range("A20").Value = -0.063439365214835952
$ws.Range("B20").Value = 0.063262066135066242
$ws.Range("A21").Value = -0.0090044366559456357
$ws.Range("B21").Value = 0.0089999998411296467
$ws.Range("A22").Value = -0.093951553909791485
$ws.Range("B22").Value = 0.093636870780676063
$ws.Range("A23").Value = -0.084636870939632125
$ws.Range("B23").Value = 0.084127414340007256
$ws.Range("A24").Value = -0.042127414579925748
$ws.Range("B24").Value = 0.041999999758715489
$ws.Range("A25").Value = -0.094986889705857891
$ws.Range("B25").Value = 0.094735788406676846
$ws.Range("A26").Value = -0.088735788567870344
$ws.Range("B26").Value = 0.088413403179487204
$ws.Range("A27").Value = -0.082413403341683455
$ws.Range("B27").Value = 0.081315187857679216
$ws.Range("A28").Value = -0.075315188023672874
$ws.Range("B28").Value = 0.074554077432647858
$ws.Range("A29").Value = -0.062554077615535775
$ws.Range("B29").Value = 0.062176068958368802
$ws.Range("A30").Value = -0.042176069161545371
$ws.Range("B30").Value = 0.042021293478979427
$ws.Range("A31").Value = -0.027021293672179425
$ws.Range("B31").Value = 0.027001055806827878
$ws.Range("A32").Value = -0.006001056014597772
$ws.Range("B32").Value = 0.0059999998276758504